$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.121.69"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.668.77"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "'210.42"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "'0.5215"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'0.2618"
$ws.Range("E8").Value = "  -2.22%  "
$ws.Range("D9").Value = "'0.06323"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'21.16"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").Value = "'0.07542"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "1.673.99"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "'4.425"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "'0.5458"
$ws.Range("E14").Value = "  -4.36%  "
$ws.Range("D15").Value = "'0.000008019"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").Value = "'66.46"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("D17").Value = "26.163.81"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'4.746"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").Value = "'187.49"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").Value = "'10.30"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").Value = "'6.243"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "'149.74"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("D25").Value = "'0.1239"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").Value = "'7.487"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").Value = "'15.77"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "'0.06293"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "'1.355"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'3.511"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'3.422"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").Value = "'1.649"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").Value = "'1.003"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "'0.6013"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "'2.767"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "'2.396"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "1.116.30"
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").Value = "'0.01614"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "'6.061"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").Value = "'0.8626"
$ws.Range("E41").Value = "  -1.44%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "'100.53"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "1.823.18"
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("D46").Value = "'55.54"
$ws.Range("E46").Value = "  -3.05%  "
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "'8.056"
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("D49").Value = "'0.05254"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "'0.4240"
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "'5.915"
$ws.Range("E51").Value = "  -0.89%  "
